$wb = $excel.ActiveWorkbook

# Insert a brand-new worksheet before the current first sheet (Personnel),
# so it becomes the new first sheet of the workbook.
$firstSheet = $wb.Worksheets.Item(1)
$newSheet = $wb.Worksheets.Add($firstSheet)
$newSheet.Name = "CategoricalVariables"

# Header row
$newSheet.Range("A1").Value = "attributeName"
$newSheet.Range("B1").Value = "code"
$newSheet.Range("C1").Value = "definition"

# Data rows
$newSheet.Range("A2").Value = "toi_source"
$newSheet.Range("B2").Value = "toi_niskin"
$newSheet.Range("C2").Value = "sample bottle was filled from a Niskin bottle on CTD rosette"

$newSheet.Range("A3").Value = "toi_source"
$newSheet.Range("B3").Value = "toi_underway"
$newSheet.Range("C3").Value = "sample bottle was filled from the ship's underway system"

# Make the new sheet the active/selected one, matching the tabSelected view
# state moving to the first sheet, with a selection range of A1:E5.
$newSheet.Range("A1:E5").Select() | Out-Null
$newSheet.Activate() | Out-Null
